# Apply recalculated "Tabla_Victorias_DM_CORREGIDO" results.
# Updates the B (Modelo), C (Score_Neto), D (Porcentaje_Victorias),
# E (Victorias) and F (Derrotas) columns on the 4 worksheets so that
# rows reflect the corrected rankings (mainly affecting LSPMW and the
# models around it).

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        [string]$modelo,
        [double]$score,
        [double]$pct,
        [double]$victorias,
        [double]$derrotas
    )
    $ws.Cells.Item($row, 2).Value = $modelo
    $ws.Cells.Item($row, 3).Value = $score
    $ws.Cells.Item($row, 4).Value = $pct
    $ws.Cells.Item($row, 5).Value = $victorias
    $ws.Cells.Item($row, 6).Value = $derrotas
}

# --- Sheet "General" ---
$ws1 = $wb.Worksheets.Item("General")
Set-Row $ws1 4 "AV-MCPS" 2    50   4 2
Set-Row $ws1 5 "DeepAR"  2    50   4 2
Set-Row $ws1 6 "MCPS"    2    50   4 2
Set-Row $ws1 7 "LSPMW"   -2   37.5 3 5

# --- Sheet "Est_Lin_ARMA" ---
$ws2 = $wb.Worksheets.Item("Est_Lin_ARMA")
Set-Row $ws2 2 "Sieve Bootstrap"     8  100  8 0
Set-Row $ws2 3 "DeepAR"              6  87.5 7 1
Set-Row $ws2 5 "MCPS"                2  50   4 2
Set-Row $ws2 6 "LSPM"                1  50   4 3
Set-Row $ws2 7 "Block Bootstrapping" -2 37.5 3 5
Set-Row $ws2 8 "AREPD"               -4 25   2 6
Set-Row $ws2 9 "EnCQR-LSTM"          -7 0    0 7
Set-Row $ws2 10 "LSPMW"              -7 0    0 7

# --- Sheet "Est_NoLin_SETAR" ---
$ws3 = $wb.Worksheets.Item("Est_NoLin_SETAR")
Set-Row $ws3 2 "DeepAR"              6  75   6 0
Set-Row $ws3 3 "Sieve Bootstrap"     6  75   6 0
Set-Row $ws3 4 "Block Bootstrapping" 5  62.5 5 0
Set-Row $ws3 5 "AV-MCPS"             -1 12.5 1 2
Set-Row $ws3 6 "AREPD"               -2 12.5 1 3
Set-Row $ws3 7 "LSPM"                -2 12.5 1 3
Set-Row $ws3 8 "LSPMW"               -2 12.5 1 3
Set-Row $ws3 9 "MCPS"                -2 12.5 1 3

# --- Sheet "NoEst_Lin_ARIMA" ---
$ws4 = $wb.Worksheets.Item("NoEst_Lin_ARIMA")
Set-Row $ws4 3 "LSPM"    6  87.5 7 1
Set-Row $ws4 4 "AV-MCPS" 2  50   4 2
Set-Row $ws4 5 "MCPS"    2  50   4 2
Set-Row $ws4 6 "DeepAR"  1  37.5 3 2
Set-Row $ws4 7 "LSPMW"   -1 37.5 3 4

$wb.Save()
